$wb = $excel.ActiveWorkbook

# --- open_buy_orders sheet: remove the two data rows (rows 2 and 3) ---
$buyWs = $wb.Worksheets.Item("open_buy_orders")
$buyWs.Rows.Item(3).Delete()
$buyWs.Rows.Item(2).Delete()

# --- open_sell_orders sheet: update txid in A2 ---
$sellWs = $wb.Worksheets.Item("open_sell_orders")
$sellWs.Range("A2").Value = "OEIQ7F-ABHMY-7A6R26"
